$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.961.80'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.562.09'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.57'
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.06'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0599'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '1.784.90'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.564.73'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.09'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '26.958.60'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.56'
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.22'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.57'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("E32").Value = '  -0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("D34").Value = '1.420.76'
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("E35").Value = '  +11.39%  '
$ws.Range("E36").Value = '  +2.57%  '
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.533'
$ws.Range("E39").Value = '  +1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.79'
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.808'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("E44").Value = '  +1.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.84'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = '1.698.26'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.36'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("E49").Value = '  +2.01%  '
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  -0.71%  '
